# Auto-generated edit script: updates cached market-price / profit values
# across the Leve profit tracking sheets (scheduled data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 733
$ws.Range("J17").Value = 733
$ws.Range("L17").Value = 2199
$ws.Range("N17").Value = -2535
$ws.Range("H19").Value = 2208.25
$ws.Range("I19").Value = 2085.8
$ws.Range("J19").Value = 2330.7
$ws.Range("K19").Value = 2085.8
$ws.Range("L19").Value = 2330.7
$ws.Range("M19").Value = -1910.8
$ws.Range("N19").Value = -2680.7
$ws.Range("H127").Value = 1021.5
$ws.Range("I127").Value = 1053.4286
$ws.Range("K127").Value = 3160.2858
$ws.Range("M127").Value = 1799.7142
$ws.Range("H132").Value = 75305.5
$ws.Range("I132").Value = 80713.62
$ws.Range("K132").Value = 242140.86
$ws.Range("M132").Value = -239610.86
$ws.Range("H137").Value = 970
$ws.Range("I137").Value = 970
$ws.Range("K137").Value = 2910
$ws.Range("M137").Value = -360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9260235
$ws.Range("I74").Value = 18518518
$ws.Range("J74").Value = 6945664
$ws.Range("K74").Value = 18518518
$ws.Range("L74").Value = 6945664
$ws.Range("M74").Value = -18517644
$ws.Range("N74").Value = -6947412
$ws.Range("H77").Value = 9260235
$ws.Range("I77").Value = 18518518
$ws.Range("J77").Value = 6945664
$ws.Range("K77").Value = 92592590
$ws.Range("L77").Value = 34728320
$ws.Range("M77").Value = -92588222
$ws.Range("N77").Value = -34737056
$ws.Range("H97").Value = 1424.6
$ws.Range("J97").Value = 2419.75
$ws.Range("L97").Value = 2419.75
$ws.Range("N97").Value = -3411.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 43437.668
$ws.Range("J76").Value = 43437.668
$ws.Range("L76").Value = 43437.668
$ws.Range("N76").Value = -44067.668
$ws.Range("H79").Value = 43437.668
$ws.Range("J79").Value = 43437.668
$ws.Range("L79").Value = 43437.668
$ws.Range("N79").Value = -45621.668
$ws.Range("H86").Value = 2682.5557
$ws.Range("I86").Value = 2106.1428
$ws.Range("K86").Value = 2106.1428
$ws.Range("M86").Value = -983.1428000000001
$ws.Range("H89").Value = 2682.5557
$ws.Range("I89").Value = 2106.1428
$ws.Range("K89").Value = 10530.714
$ws.Range("M89").Value = -4914.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 6401.6
$ws.Range("I2").Value = 7727
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 7727
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = -7614
$ws.Range("N2").Value = -1326
$ws.Range("H3").Value = 3646.3333
$ws.Range("I3").Value = 4944
$ws.Range("K3").Value = 4944
$ws.Range("M3").Value = -4831
$ws.Range("H7").Value = 498.5
$ws.Range("I7").Value = 143.08333
$ws.Range("J7").Value = 2631
$ws.Range("K7").Value = 143.08333
$ws.Range("L7").Value = 2631
$ws.Range("M7").Value = -30.08332999999999
$ws.Range("N7").Value = -2857
$ws.Range("H22").Value = 20333.666
$ws.Range("I22").Value = 25499.5
$ws.Range("K22").Value = 25499.5
$ws.Range("M22").Value = -25149.5
$ws.Range("H23").Value = 33336.668
$ws.Range("I23").Value = 30000
$ws.Range("J23").Value = 35005
$ws.Range("K23").Value = 30000
$ws.Range("L23").Value = 35005
$ws.Range("M23").Value = -29760
$ws.Range("N23").Value = -35485
$ws.Range("H27").Value = 33336.668
$ws.Range("I27").Value = 30000
$ws.Range("J27").Value = 35005
$ws.Range("K27").Value = 30000
$ws.Range("L27").Value = 35005
$ws.Range("M27").Value = -29808
$ws.Range("N27").Value = -35389
$ws.Range("H58").Value = 2042.5454
$ws.Range("I58").Value = 2265.7827
$ws.Range("K58").Value = 2265.7827
$ws.Range("M58").Value = -2062.7827
$ws.Range("H132").Value = 2921.5715
$ws.Range("I132").Value = 2551.3333
$ws.Range("J132").Value = 3199.25
$ws.Range("K132").Value = 7653.999899999999
$ws.Range("L132").Value = 9597.75
$ws.Range("M132").Value = -5123.999899999999
$ws.Range("N132").Value = -14657.75
$ws.Range("H134").Value = 3139.75
$ws.Range("I134").Value = 2872.2222
$ws.Range("J134").Value = 3942.3333
$ws.Range("K134").Value = 8616.6666
$ws.Range("L134").Value = 11826.9999
$ws.Range("M134").Value = -6081.6666
$ws.Range("N134").Value = -16896.9999
$ws.Range("H136").Value = 2042.5454
$ws.Range("I136").Value = 2265.7827
$ws.Range("K136").Value = 6797.348100000001
$ws.Range("M136").Value = -4247.348100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 7749.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7749.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 23248.5
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -23430.5
$ws.Range("H80").Value = 4467.5
$ws.Range("I80").Value = 4298.5
$ws.Range("K80").Value = 12895.5
$ws.Range("M80").Value = -11959.5
$ws.Range("H83").Value = 4467.5
$ws.Range("I83").Value = 4298.5
$ws.Range("K83").Value = 38686.5
$ws.Range("M83").Value = -34006.5
$ws.Range("H108").Value = 2899.5
$ws.Range("I108").Value = 2899.5
$ws.Range("K108").Value = 8698.5
$ws.Range("M108").Value = -5818.5
$ws.Range("H120").Value = 20194
$ws.Range("I120").Value = 20194
$ws.Range("K120").Value = 60582
$ws.Range("M120").Value = -55744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3130.5833
$ws.Range("I102").Value = 2292.25
$ws.Range("J102").Value = 3549.75
$ws.Range("K102").Value = 2292.25
$ws.Range("L102").Value = 3549.75
$ws.Range("M102").Value = -670.25
$ws.Range("N102").Value = -6793.75
$ws.Range("H132").Value = 2511.7144
$ws.Range("I132").Value = 2187.7144
$ws.Range("K132").Value = 6563.1432
$ws.Range("M132").Value = -4033.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5653.4546
$ws.Range("I9").Value = 877.25
$ws.Range("J9").Value = 8382.714
$ws.Range("K9").Value = 877.25
$ws.Range("L9").Value = 8382.714
$ws.Range("M9").Value = -653.25
$ws.Range("N9").Value = -8830.714
$ws.Range("H16").Value = 3144
$ws.Range("I16").Value = 3144
$ws.Range("K16").Value = 3144
$ws.Range("M16").Value = -2974
$ws.Range("H22").Value = 10999.8
$ws.Range("I22").Value = 6665.6665
$ws.Range("J22").Value = 17501
$ws.Range("K22").Value = 6665.6665
$ws.Range("L22").Value = 17501
$ws.Range("M22").Value = -6370.6665
$ws.Range("N22").Value = -18091
$ws.Range("H27").Value = 10999.8
$ws.Range("I27").Value = 6665.6665
$ws.Range("J27").Value = 17501
$ws.Range("K27").Value = 6665.6665
$ws.Range("L27").Value = 17501
$ws.Range("M27").Value = -6558.6665
$ws.Range("N27").Value = -17715
$ws.Range("H40").Value = 3027.487
$ws.Range("I40").Value = 2582.1155
$ws.Range("J40").Value = 3918.2307
$ws.Range("K40").Value = 2582.1155
$ws.Range("L40").Value = 3918.2307
$ws.Range("M40").Value = -2446.1155
$ws.Range("N40").Value = -4190.2307
$ws.Range("H46").Value = 2075.7144
$ws.Range("J46").Value = 1754.7368
$ws.Range("L46").Value = 1754.7368
$ws.Range("N46").Value = -2130.7368
$ws.Range("H55").Value = 601
$ws.Range("I55").Value = 655.4666999999999
$ws.Range("J55").Value = 555.6111
$ws.Range("K55").Value = 655.4666999999999
$ws.Range("L55").Value = 555.6111
$ws.Range("M55").Value = -482.4666999999999
$ws.Range("N55").Value = -901.6111
$ws.Range("H98").Value = 66152.5
$ws.Range("J98").Value = 66152.5
$ws.Range("L98").Value = 66152.5
$ws.Range("N98").Value = -72142.5
$ws.Range("H102").Value = 21894.334
$ws.Range("J102").Value = 21894.334
$ws.Range("L102").Value = 21894.334
$ws.Range("N102").Value = -28384.334
$ws.Range("H122").Value = 3856.1428
$ws.Range("I122").Value = 3799.2
$ws.Range("J122").Value = 3998.5
$ws.Range("K122").Value = 11397.6
$ws.Range("L122").Value = 11995.5
$ws.Range("M122").Value = -8947.599999999999
$ws.Range("N122").Value = -16895.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 28481.889
$ws.Range("J102").Value = 29476.715
$ws.Range("L102").Value = 29476.715
$ws.Range("N102").Value = -35966.715
